# Automatische test-sync: 2025-07-31 21:50:50
# Adds a new test-mail log entry to the "Logs" sheet and bumps the
# "Productinformatie" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$newRow = 15

$values = @(
    "Kun je mij de datasheet van de VentiQ-250 sturen?",
    "mailmind.test@zohomail.eu",
    "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?",
    "Productinformatie",
    "Beste klant,`nBedankt voor je interesse in de VentiQ-250. Helaas kunnen we op basis van je e-mailadres geen datasheet vinden. Zou je ons kunnen voorzien van meer informatie, zoals je volledige naam, bedrijfsnaam of eventuele andere gegevens waaronder de datasheet geregistreerd staat? Hiermee kunnen we je beter van dienst zijn en de datasheet naar je opsturen.`nMet vriendelijke groet,`n[Jouw naam]`nE-mailassistent",
    "2025-07-31 21:50:21",
    "Ja",
    "Nee",
    "Ja",
    "Nee"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $logs.Cells.Item($newRow, $col).Value = $values[$i]
}

# The multi-line reply text makes the engine auto-fit the row height;
# put it back to the sheet's default (no explicit row height), matching
# every other data row.
$logs.Rows.Item($newRow).AutoFit()

# Extend the conditional formatting ranges so they keep covering the
# whole data range, now through row 15 (one ModifyAppliesToRange call per
# block is enough to re-point every rule sharing that sqref).
$logs.Range("D2:D14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D15"))
$logs.Range("G2:G14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G15"))
$logs.Range("H2:H14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H15"))
$logs.Range("I2:I14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I15"))
$logs.Range("J2:J14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J15"))

# Keep the category tally on the Dashboard sheet in sync with the new row.
$dash.Range("B3").Value = 4
